{"js": "// Update the date paragraph (first paragraph in the body) in place,\n// preserving its run formatting (font/size/alignment).\nconst firstParagraphs = context.document.body.paragraphs;\nfirstParagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = firstParagraphs.items[0];\ntitlePara.insertText(\"2025-09-25 Thursday\", Word.InsertLocation.replace);\n\n// Update the five \"answer\" rows of the table (row indices 0, 4, 8, 12, 16 -\n// the other rows are blank work rows). Writing to TableCell.value replaces\n// the cell's text while keeping the existing cell/run formatting intact.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRowValues = {\n  0: [\"21\u00f78=2, 5\", \"20\u00f75=4, 0\", \"69\u00f74=17, 1\", \"81\u00f78=10, 1\", \"74\u00f77=10, 4\"],\n  4: [\"16\u00f74=4, 0\", \"87\u00f74=21, 3\", \"15\u00f73=5, 0\", \"56\u00f75=11, 1\", \"69\u00f77=9, 6\"],\n  8: [\"12\u00f73=4, 0\", \"67\u00f74=16, 3\", \"53\u00f78=6, 5\", \"66\u00f73=22, 0\", \"33\u00f77=4, 5\"],\n  12: [\"25\u00f75=5, 0\", \"70\u00f73=23, 1\", \"31\u00f75=6, 1\", \"38\u00f78=4, 6\", \"62\u00f73=20, 2\"],\n  16: [\"76\u00f74=19, 0\", \"91\u00f72=45, 1\", \"24\u00f76=4, 0\", \"92\u00f78=11, 4\", \"84\u00f76=14, 0\"],\n};\n\nfor (const rowIndex of Object.keys(newRowValues)) {\n  const rowIdx = Number(rowIndex);\n  const values = newRowValues[rowIndex];\n  for (let colIdx = 0; colIdx < values.length; colIdx++) {\n    const cell = table.getCell(rowIdx, colIdx);\n    cell.value = values[colIdx];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line (first paragraph), preserving its run formatting.\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Text = \"2025-09-25 Thursday\"\n\n# Update the five \"answer\" rows of the table (Word COM tables/cells are\n# 1-based: rows 1, 5, 9, 13, 17 hold data; the rows in between are blank\n# work rows left untouched). Setting Cell.Range.Text replaces only the\n# text, keeping the existing cell/run formatting intact.\n$table = $d.Tables.Item(1)\n\n$newRowValues = @{\n    1  = @(\"21\u00f78=2, 5\", \"20\u00f75=4, 0\", \"69\u00f74=17, 1\", \"81\u00f78=10, 1\", \"74\u00f77=10, 4\")\n    5  = @(\"16\u00f74=4, 0\", \"87\u00f74=21, 3\", \"15\u00f73=5, 0\", \"56\u00f75=11, 1\", \"69\u00f77=9, 6\")\n    9  = @(\"12\u00f73=4, 0\", \"67\u00f74=16, 3\", \"53\u00f78=6, 5\", \"66\u00f73=22, 0\", \"33\u00f77=4, 5\")\n    13 = @(\"25\u00f75=5, 0\", \"70\u00f73=23, 1\", \"31\u00f75=6, 1\", \"38\u00f78=4, 6\", \"62\u00f73=20, 2\")\n    17 = @(\"76\u00f74=19, 0\", \"91\u00f72=45, 1\", \"24\u00f76=4, 0\", \"92\u00f78=11, 4\", \"84\u00f76=14, 0\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $values = $newRowValues[$rowIndex]\n    for ($colIndex = 1; $colIndex -le $values.Length; $colIndex++) {\n        $cell = $table.Cell($rowIndex, $colIndex)\n        $cell.Range.Text = $values[$colIndex - 1]\n    }\n}\n"}
